# #578, remove YYYYMMDD Token, add predicate
#
# The "TokenTypes" sheet is a generated enum table: column A is a category
# marker, column B is an auto-incrementing id (=prev+1), column C is the
# token name (shared string) and column D rebuilds the C# enum line from
# B/C. Row 442 held the YYYYMMDD token; it needs to disappear and every
# row below it (443..454) needs to shift up by one, which (because B/D are
# formulas) happens automatically once the token names in column C are
# moved up - Excel recalculates the ids and rebuilt enum text for us, and
# drops "YYYYMMDD" from the shared-string table since nothing references
# it any more.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Token names currently in C443:C454 (i.e. everything below YYYYMMDD),
# shifted up into C442:C453 - this both removes the YYYYMMDD row and closes
# the gap left behind.
$tokens = @(
    "TYPEDEF",
    "STRONG",
    "UNSAFE",
    "PUBLIC",
    "PRIVATE",
    "IN_OUT",
    "STRICT",
    "QuestionMark",
    "CompilerDirective",
    "CopyImportDirective",
    "ReplaceDirective",
    "ContinuationTokenGroup"
)

for ($i = 0; $i -lt $tokens.Length; $i++) {
    $row = 442 + $i
    $ws.Range("C$row").Value = $tokens[$i]
}

# Row 442 (now TYPEDEF) is a plain token row like 443-454 already were, so
# it loses the "A" category marker that only decorated the YYYYMMDD row.
$ws.Range("A442").ClearContents()

# The last row (454, ContinuationTokenGroup) is now a duplicate of the row
# that got shifted into 453 - delete it outright to shrink the table back
# to A1:D453.
$ws.Rows("454:454").Delete()

# Park the selection where the author's Excel session ended up after the
# edit.
$ws.Range("C455").Select() | Out-Null
